$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the header text in A1 and B1
$a1 = $ws.Range("A1").Value2
$b1 = $ws.Range("B1").Value2
$ws.Range("A1").Value = $b1
$ws.Range("B1").Value = $a1

# Update the active selection to F12 (cosmetic, matches the recorded edit)
$ws.Activate()
$ws.Range("F12").Select()
